$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# "Create table merge column": duplicate the existing row and vertically
# merge the first column ("milad") across both rows, while the second
# column ("reza") stays as two independent cells.

# 1) Add a second row with the same cell text as row 1.
$table.Rows.Add()
$table.Cell(2, 1).Range.Text = "milad"
$table.Cell(2, 2).Range.Text = "reza"

# 2) Vertically merge the first-column cells of row 1 and row 2: the row-1
#    cell becomes the merge anchor (vMerge val="restart") and the row-2
#    cell becomes the merge continuation (vMerge).
$table.Cell(1, 1).Merge($table.Cell(2, 1))

# Cell.Merge clears the continuation cell's paragraph text as a side
# effect; put "milad" back so both rows of the merged column read "milad".
$table.Rows.Item(2).Cells.Item(1).Range.Text = "milad"

# Re-apply the table as clean WordprocessingML so the merged row is a
# plain <w:tr> (no leftover rsid/paraId bookkeeping from the row-insert
# operation), matching a freshly authored table.
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"
$tableXml = "<w:tbl xmlns:w='$w'>" +
  "<w:tblPr>" +
    "<w:tblStyle w:val='TableGrid'/>" +
    "<w:tblW w:w='0' w:type='auto'/>" +
    "<w:tblLook w:val='04A0' w:firstRow='1' w:lastRow='0' w:firstColumn='1' w:lastColumn='0' w:noHBand='0' w:noVBand='1'/>" +
  "</w:tblPr>" +
  "<w:tblGrid>" +
    "<w:gridCol w:w='4675'/>" +
    "<w:gridCol w:w='4675'/>" +
  "</w:tblGrid>" +
  "<w:tr>" +
    "<w:tc>" +
      "<w:tcPr><w:tcW w:w='4675' w:type='dxa'/><w:vMerge w:val='restart'/></w:tcPr>" +
      "<w:p><w:r><w:t>milad</w:t></w:r></w:p>" +
    "</w:tc>" +
    "<w:tc>" +
      "<w:tcPr><w:tcW w:w='4675' w:type='dxa'/></w:tcPr>" +
      "<w:p><w:r><w:t>reza</w:t></w:r></w:p>" +
    "</w:tc>" +
  "</w:tr>" +
  "<w:tr>" +
    "<w:tc>" +
      "<w:tcPr><w:tcW w:w='4675' w:type='dxa'/><w:vMerge/></w:tcPr>" +
      "<w:p><w:r><w:t>milad</w:t></w:r></w:p>" +
    "</w:tc>" +
    "<w:tc>" +
      "<w:tcPr><w:tcW w:w='4675' w:type='dxa'/></w:tcPr>" +
      "<w:p><w:r><w:t>reza</w:t></w:r></w:p>" +
    "</w:tc>" +
  "</w:tr>" +
"</w:tbl>"

$d.Content.InsertXML($tableXml) | Out-Null
